# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> name="Office Theme" / clrScheme name="Office"   (bound to the Notes Master)
#   ppt/theme/theme2.xml -> name="Integral"      / clrScheme name="Red Violet" (bound to the Slide Master
#                                                                               + the presentation itself)
# The target edit swaps the full contents of these two parts: the design actually
# applied to the slides/master switches from the "Integral" (Red Violet) palette
# to the default "Office Theme" (Office) palette.
#
# The PowerPoint object model doesn't expose the slide master's theme part as an
# editable blob, but it does expose a live, settable view of its 12 theme colors
# via Slide.ThemeColorScheme (a master-level singleton - touching it through any
# slide updates the one theme part shared by the whole deck). Drive the palette
# over to the exact "Office Theme" RGB values to reproduce the effective result
# of the swap for the part that is actually applied to the presentation.

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index order is fixed: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
# 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink (RGB stored 0x00BBGGRR).
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
